# Atualizações do plano de acao
#
# 1) "Novo Banco de Dados com DER - Cristhian e Kaue"
#       -> "DER" + " do" + " Banco de Dados - Cristhian e Kaue"   (3 runs)
# 2) "Script" + "- Larissa e Tabata"
#       -> "Script" + " Banco de Dados" + "- Larissa e Tabata"    (3 runs)

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: collapse "Novo Banco de Dados com DER" down to just "DER",
# then re-insert " do" and " Banco de Dados" as their own runs (editing
# right-to-left so earlier inserts are not disturbed by later ones).
# ---------------------------------------------------------------------

$rngDer = $d.Content
$rngDer.Find.Execute("Novo Banco de Dados com DER", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "DER", 2) | Out-Null

# Insert " Banco de Dados" right after "DER" (before " - Cristhian e Kaue")
$rngFull = $d.Content
$rngFull.Find.Execute("DER - Cristhian e Kaue") | Out-Null
$afterDer = $rngFull.Start + 3
$insPoint1 = $d.Range($afterDer, $afterDer)
$insPoint1.InsertBefore(" Banco de Dados")
$insPoint1.Font.Bold = 1
$insPoint1.Font.Bold = 0

# Insert " do" right after "DER" (now immediately before the run we just added)
$rngDerOnly = $d.Content
$rngDerOnly.Find.Execute("DER") | Out-Null
$derEnd = $rngDerOnly.End
$insPoint2 = $d.Range($derEnd, $derEnd)
$insPoint2.InsertBefore(" do")
$insPoint2.Font.Bold = 1
$insPoint2.Font.Bold = 0

# ---------------------------------------------------------------------
# Change 2: insert " Banco de Dados" right after "Script" and before
# "- Larissa e Tabata" as its own run.
# ---------------------------------------------------------------------

$rngScript = $d.Content
$rngScript.Find.Execute("Script") | Out-Null
$scriptEnd = $rngScript.End
$insPoint3 = $d.Range($scriptEnd, $scriptEnd)
$insPoint3.InsertBefore(" Banco de Dados")
$insPoint3.Font.Bold = 1
$insPoint3.Font.Bold = 0

Write-Output $d.Content.Text
